# Config.xlsx update:
#  - Insert two new Settings rows (CustomTransactionNumbers, OutputDataFilePath)
#    above the existing "wcm_CredentialsName" row.
#  - Give the Value column of the new rows a Text number format so values
#    like "2-8" are stored verbatim.
#  - Move the active selection to C5 (the description cell of the new
#    OutputDataFilePath row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Make room for the two new rows; existing rows 4.. shift down to 6..
$ws.Rows("4:5").Insert()

# Row 4: CustomTransactionNumbers
$ws.Range("A4").Value = "CustomTransactionNumbers"
$ws.Range("C4").Value = "Empty, if all rows of Input File shall be processed. A semicolon separated list of Excel Row Numbers (3;6;20) or a hypen separated range of Excel Row Numbers (3-20), if only these Excel Rows shall be processed."
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2-8"
$ws.Rows(4).RowHeight = 45

# Row 5: OutputDataFilePath
$ws.Range("A5").Value = "OutputDataFilePath"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "Data\Output\"
$ws.Range("C5").Value = "Path, where the Output Excel Workbooks are saved"

# Match the author's final selection on the sheet.
$ws.Activate()
$ws.Range("C5").Select()
